# Update the arithmetic problems in the worksheet table.
# Each old value is unique within the document, so a simple
# Find/Replace (wdReplaceAll) per pair is safe and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("39÷9=", "66÷9="),
    @("13÷5=", "30÷6="),
    @("91÷2=", "99÷7="),
    @("49÷2=", "82÷7="),
    @("23÷8=", "85÷5="),
    @("85÷9=", "87÷5="),
    @("47÷2=", "30÷2="),
    @("31÷9=", "42÷2="),
    @("41÷2=", "66÷3="),
    @("14÷8=", "71÷3="),
    @("63÷9=", "52÷6="),
    @("69÷9=", "53÷6="),
    @("89÷4=", "52÷9="),
    @("94÷5=", "37÷8="),
    @("21÷3=", "19÷9="),
    @("15÷2=", "43÷7="),
    @("68÷2=", "12÷7="),
    @("53÷7=", "86÷6="),
    @("75÷2=", "26÷4="),
    @("85÷4=", "66÷7="),
    @("62÷6=", "67÷2="),
    @("97÷2=", "60÷2="),
    @("85÷7=", "65÷7="),
    @("68÷7=", "81÷9="),
    @("59÷3=", "80÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
